$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "Information Technology" = "Technology"
    "Consumer Discretionary" = "Consumer Cyclical"
    "Financials"             = "Financial Services"
    "Materials"              = "Basic Materials"
    "Communication"          = "Communication Services"
}

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
